$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 11477
$ws.Range("I2").Value = 1856.4286
$ws.Range("K2").Value = 1856.4286
$ws.Range("M2").Value = -1743.4286
$ws.Range("H19").Value = 4876.25
$ws.Range("I19").Value = 5756.5
$ws.Range("K19").Value = 5756.5
$ws.Range("M19").Value = -5581.5
$ws.Range("H70").Value = 5054.5
$ws.Range("I70").Value = 4559.625
$ws.Range("J70").Value = 5337.2856
$ws.Range("K70").Value = 13678.875
$ws.Range("L70").Value = 16011.8568
$ws.Range("M70").Value = -13408.875
$ws.Range("N70").Value = -16551.8568
$ws.Range("H73").Value = 5054.5
$ws.Range("I73").Value = 4559.625
$ws.Range("J73").Value = 5337.2856
$ws.Range("K73").Value = 13678.875
$ws.Range("L73").Value = 16011.8568
$ws.Range("M73").Value = -12742.875
$ws.Range("N73").Value = -17883.8568
$ws.Range("H100").Value = 52962.2
$ws.Range("I100").Value = 72232.28999999999
$ws.Range("K100").Value = 72232.28999999999
$ws.Range("M100").Value = -71691.28999999999
$ws.Range("H111").Value = 3335434
$ws.Range("I111").Value = 3180.125
$ws.Range("J111").Value = 4547162.5
$ws.Range("K111").Value = 9540.375
$ws.Range("L111").Value = 13641487.5
$ws.Range("M111").Value = -6473.375
$ws.Range("N111").Value = -13647621.5
$ws.Range("H112").Value = 2848.3794
$ws.Range("I112").Value = 1696.3334
$ws.Range("J112").Value = 2981.3076
$ws.Range("K112").Value = 5089.0002
$ws.Range("L112").Value = 8943.9228
$ws.Range("M112").Value = -3981.0002
$ws.Range("N112").Value = -11159.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 2263.25
$ws.Range("J12").Value = 6004
$ws.Range("L12").Value = 6004
$ws.Range("N12").Value = -6350
$ws.Range("H61").Value = 5784.091
$ws.Range("I61").Value = 4541.722
$ws.Range("K61").Value = 4541.722
$ws.Range("M61").Value = -4329.722
$ws.Range("H74").Value = 5478.5293
$ws.Range("I74").Value = 3155.75
$ws.Range("K74").Value = 3155.75
$ws.Range("M74").Value = -2281.75
$ws.Range("H77").Value = 5478.5293
$ws.Range("I77").Value = 3155.75
$ws.Range("K77").Value = 15778.75
$ws.Range("M77").Value = -11410.75
$ws.Range("H102").Value = 2222.12
$ws.Range("I102").Value = 2222.12
$ws.Range("K102").Value = 2222.12
$ws.Range("M102").Value = -600.1199999999999
$ws.Range("H122").Value = 1671.1428
$ws.Range("I122").Value = 1449.6666
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4348.9998
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1898.9998
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 3267.075
$ws.Range("I132").Value = 2867.8
$ws.Range("J132").Value = 3932.5334
$ws.Range("K132").Value = 8603.400000000001
$ws.Range("L132").Value = 11797.6002
$ws.Range("M132").Value = -6073.400000000001
$ws.Range("N132").Value = -16857.6002
$ws.Range("H136").Value = 5784.091
$ws.Range("I136").Value = 4541.722
$ws.Range("K136").Value = 13625.166
$ws.Range("M136").Value = -11075.166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 3054.111
$ws.Range("J64").Value = 3061.375
$ws.Range("L64").Value = 3061.375
$ws.Range("N64").Value = -3511.375
$ws.Range("H67").Value = 3054.111
$ws.Range("J67").Value = 3061.375
$ws.Range("L67").Value = 3061.375
$ws.Range("N67").Value = -4621.375
$ws.Range("H94").Value = 2812.2666
$ws.Range("I94").Value = 2242.6667
$ws.Range("J94").Value = 3666.6667
$ws.Range("K94").Value = 2242.6667
$ws.Range("L94").Value = 3666.6667
$ws.Range("M94").Value = -1791.6667
$ws.Range("N94").Value = -4568.6667
$ws.Range("H107").Value = 3091.9167
$ws.Range("J107").Value = 4987.25
$ws.Range("L107").Value = 4987.25
$ws.Range("N107").Value = -8827.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H44").Value = 15469
$ws.Range("I44").Value = 13815
$ws.Range("J44").Value = 17950
$ws.Range("K44").Value = 13815
$ws.Range("L44").Value = 17950
$ws.Range("M44").Value = -13373
$ws.Range("N44").Value = -18834
$ws.Range("H53").Value = 49561.332
$ws.Range("J53").Value = 49561.332
$ws.Range("L53").Value = 49561.332
$ws.Range("N53").Value = -50775.332
$ws.Range("H86").Value = 9939.375
$ws.Range("I86").Value = 9101.4
$ws.Range("J86").Value = 11336
$ws.Range("K86").Value = 9101.4
$ws.Range("L86").Value = 11336
$ws.Range("M86").Value = -7978.4
$ws.Range("N86").Value = -13582
$ws.Range("H89").Value = 9939.375
$ws.Range("I89").Value = 9101.4
$ws.Range("J89").Value = 11336
$ws.Range("K89").Value = 45507
$ws.Range("L89").Value = 56680
$ws.Range("M89").Value = -39891
$ws.Range("N89").Value = -67912
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
$ws.Range("H105").Value = 2179.4
$ws.Range("I105").Value = 1974.25
$ws.Range("K105").Value = 1974.25
$ws.Range("M105").Value = -227.25
$ws.Range("H107").Value = 1805.4445
$ws.Range("I107").Value = 1299.8667
$ws.Range("J107").Value = 4333.3335
$ws.Range("K107").Value = 1299.8667
$ws.Range("L107").Value = 4333.3335
$ws.Range("M107").Value = 620.1333
$ws.Range("N107").Value = -8173.3335
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H134").Value = 13825.833
$ws.Range("I134").Value = 14017
$ws.Range("J134").Value = 13328.8
$ws.Range("K134").Value = 42051
$ws.Range("L134").Value = 39986.39999999999
$ws.Range("M134").Value = -39516
$ws.Range("N134").Value = -45056.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 901
$ws.Range("I8").Value = 901
$ws.Range("K8").Value = 2703
$ws.Range("M8").Value = -2564
$ws.Range("H80").Value = 1520.6666
$ws.Range("I80").Value = 1677.6
$ws.Range("J80").Value = 1324.5
$ws.Range("K80").Value = 5032.799999999999
$ws.Range("L80").Value = 3973.5
$ws.Range("M80").Value = -4096.799999999999
$ws.Range("N80").Value = -5845.5
$ws.Range("H83").Value = 1520.6666
$ws.Range("I83").Value = 1677.6
$ws.Range("J83").Value = 1324.5
$ws.Range("K83").Value = 15098.4
$ws.Range("L83").Value = 11920.5
$ws.Range("M83").Value = -10418.4
$ws.Range("N83").Value = -21280.5
$ws.Range("H99").Value = 500
$ws.Range("I99").Value = 500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 746
$ws.Range("N99").ClearContents()
$ws.Range("H129").Value = 1429.3334
$ws.Range("I129").Value = 650.1818
$ws.Range("K129").Value = 1950.5454
$ws.Range("M129").Value = 3049.4546
$ws.Range("H132").Value = 1466.6666
$ws.Range("I132").Value = 1466.6666
$ws.Range("K132").Value = 13199.9994
$ws.Range("M132").Value = -10669.9994
$ws.Range("H138").Value = 5558301
$ws.Range("I138").Value = 10001442
$ws.Range("K138").Value = 30004326
$ws.Range("M138").Value = -29999186

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 49927.332
$ws.Range("J49").Value = 49927.332
$ws.Range("L49").Value = 49927.332
$ws.Range("N49").Value = -50295.332
$ws.Range("H97").Value = 741.3182
$ws.Range("I97").Value = 762.5625
$ws.Range("J97").Value = 684.6667
$ws.Range("K97").Value = 762.5625
$ws.Range("L97").Value = 684.6667
$ws.Range("M97").Value = -266.5625
$ws.Range("N97").Value = -1676.6667
$ws.Range("H107").Value = 199.33333
$ws.Range("I107").Value = 199
$ws.Range("K107").Value = 199
$ws.Range("M107").Value = 1721

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 113270.1
$ws.Range("I7").Value = 150386.58
$ws.Range("K7").Value = 150386.58
$ws.Range("M7").Value = -150274.58
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H126").Value = 113270.1
$ws.Range("I126").Value = 150386.58
$ws.Range("K126").Value = 451159.74
$ws.Range("M126").Value = -448689.74

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 31499.25
$ws.Range("I34").Value = 33998
$ws.Range("J34").Value = 30666.334
$ws.Range("K34").Value = 33998
$ws.Range("L34").Value = 30666.334
$ws.Range("M34").Value = -33795
$ws.Range("N34").Value = -31072.334
$ws.Range("H62").Value = 4067
$ws.Range("I62").Value = 3126.6
$ws.Range("J62").Value = 5634.3335
$ws.Range("K62").Value = 3126.6
$ws.Range("L62").Value = 5634.3335
$ws.Range("M62").Value = -2502.6
$ws.Range("N62").Value = -6882.3335
$ws.Range("H65").Value = 4067
$ws.Range("I65").Value = 3126.6
$ws.Range("J65").Value = 5634.3335
$ws.Range("K65").Value = 15633
$ws.Range("L65").Value = 28171.6675
$ws.Range("M65").Value = -12513
$ws.Range("N65").Value = -34411.6675
$ws.Range("H100").Value = 3134.524
$ws.Range("I100").Value = 3975.6667
$ws.Range("K100").Value = 7951.3334
$ws.Range("M100").Value = -7410.3334
$ws.Range("H132").Value = 4068.138
$ws.Range("I132").Value = 2703.8918
$ws.Range("J132").Value = 6471.8096
$ws.Range("K132").Value = 8111.6754
$ws.Range("L132").Value = 19415.4288
$ws.Range("M132").Value = -5581.6754
$ws.Range("N132").Value = -24475.4288
$ws.Range("H136").Value = 7106.643
$ws.Range("I136").Value = 2582.8333
$ws.Range("J136").Value = 10499.5
$ws.Range("K136").Value = 7748.499899999999
$ws.Range("L136").Value = 31498.5
$ws.Range("M136").Value = -5198.499899999999
$ws.Range("N136").Value = -36598.5
